$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "78+16=94"
$t.Cell(1,2).Range.Text = "14+4=18"
$t.Cell(1,3).Range.Text = "35-11=24"
$t.Cell(1,4).Range.Text = "1+5=6"
$t.Cell(1,5).Range.Text = "30+2=32"

$t.Cell(2,1).Range.Text = "67-11=56"
$t.Cell(2,2).Range.Text = "33-5=28"
$t.Cell(2,3).Range.Text = "30+0=30"
$t.Cell(2,4).Range.Text = "0+93=93"
$t.Cell(2,5).Range.Text = "45+35=80"

$t.Cell(3,1).Range.Text = "90-86=4"
$t.Cell(3,2).Range.Text = "48+10=58"
$t.Cell(3,3).Range.Text = "77+2=79"
$t.Cell(3,4).Range.Text = "48+51=99"
$t.Cell(3,5).Range.Text = "3+59=62"

$t.Cell(4,1).Range.Text = "26+58=84"
$t.Cell(4,2).Range.Text = "44-40=4"
$t.Cell(4,3).Range.Text = "41-8=33"
$t.Cell(4,4).Range.Text = "25+3=28"
$t.Cell(4,5).Range.Text = "29+66=95"

$t.Cell(5,1).Range.Text = "63+29=92"
$t.Cell(5,2).Range.Text = "95-16=79"
$t.Cell(5,3).Range.Text = "59-56=3"
$t.Cell(5,4).Range.Text = "3+9=12"
$t.Cell(5,5).Range.Text = "3+88=91"

$t.Cell(6,1).Range.Text = "11+80=91"
$t.Cell(6,2).Range.Text = "3+12=15"
$t.Cell(6,3).Range.Text = "31-19=12"
$t.Cell(6,4).Range.Text = "20-0=20"
$t.Cell(6,5).Range.Text = "13+4=17"

$t.Cell(7,1).Range.Text = "93-25=68"
$t.Cell(7,2).Range.Text = "32+16=48"
$t.Cell(7,3).Range.Text = "88-62=26"
$t.Cell(7,4).Range.Text = "83-32=51"
$t.Cell(7,5).Range.Text = "17-11=6"

$t.Cell(8,1).Range.Text = "80-62=18"
$t.Cell(8,2).Range.Text = "14+67=81"
$t.Cell(8,3).Range.Text = "45+6=51"
$t.Cell(8,4).Range.Text = "55-27=28"
$t.Cell(8,5).Range.Text = "75-51=24"

$t.Cell(9,1).Range.Text = "91-76=15"
$t.Cell(9,2).Range.Text = "74-9=65"
$t.Cell(9,3).Range.Text = "68-15=53"
$t.Cell(9,4).Range.Text = "36-32=4"
$t.Cell(9,5).Range.Text = "10+49=59"

$t.Cell(10,1).Range.Text = "28+14=42"
$t.Cell(10,2).Range.Text = "82-15=67"
$t.Cell(10,3).Range.Text = "50-33=17"
$t.Cell(10,4).Range.Text = "66+4=70"
$t.Cell(10,5).Range.Text = "34-10=24"

$t.Cell(11,1).Range.Text = "19-3=16"
$t.Cell(11,2).Range.Text = "16+45=61"
$t.Cell(11,3).Range.Text = "56-9=47"
$t.Cell(11,4).Range.Text = "51-35=16"
$t.Cell(11,5).Range.Text = "70-67=3"

$t.Cell(12,1).Range.Text = "51-8=43"
$t.Cell(12,2).Range.Text = "11+40=51"
$t.Cell(12,3).Range.Text = "44+48=92"
$t.Cell(12,4).Range.Text = "86-32=54"
$t.Cell(12,5).Range.Text = "52-17=35"

$t.Cell(13,1).Range.Text = "20+43=63"
$t.Cell(13,2).Range.Text = "62+15=77"
$t.Cell(13,3).Range.Text = "61-44=17"
$t.Cell(13,4).Range.Text = "21+28=49"
$t.Cell(13,5).Range.Text = "6+76=82"

$t.Cell(14,1).Range.Text = "52-13=39"
$t.Cell(14,2).Range.Text = "26-5=21"
$t.Cell(14,3).Range.Text = "21+35=56"
$t.Cell(14,4).Range.Text = "90-6=84"
$t.Cell(14,5).Range.Text = "55-13=42"

$t.Cell(15,1).Range.Text = "33+16=49"
$t.Cell(15,2).Range.Text = "68+17=85"
$t.Cell(15,3).Range.Text = "8+73=81"
$t.Cell(15,4).Range.Text = "5-2=3"
$t.Cell(15,5).Range.Text = "37+43=80"

$t.Cell(16,1).Range.Text = "72-5=67"
$t.Cell(16,2).Range.Text = "64-43=21"
$t.Cell(16,3).Range.Text = "85-25=60"
$t.Cell(16,4).Range.Text = "3+52=55"
$t.Cell(16,5).Range.Text = "71-45=26"

$t.Cell(17,1).Range.Text = "36+18=54"
$t.Cell(17,2).Range.Text = "10+62=72"
$t.Cell(17,3).Range.Text = "48-19=29"
$t.Cell(17,4).Range.Text = "99-79=20"
$t.Cell(17,5).Range.Text = "99-11=88"

$t.Cell(18,1).Range.Text = "79-12=67"
$t.Cell(18,2).Range.Text = "68-27=41"
$t.Cell(18,3).Range.Text = "83-41=42"
$t.Cell(18,4).Range.Text = "2+20=22"
$t.Cell(18,5).Range.Text = "8+70=78"

$t.Cell(19,1).Range.Text = "71-13=58"
$t.Cell(19,2).Range.Text = "67-45=22"
$t.Cell(19,3).Range.Text = "60-34=26"
$t.Cell(19,4).Range.Text = "81+1=82"
$t.Cell(19,5).Range.Text = "93-73=20"

$t.Cell(20,1).Range.Text = "52-16=36"
$t.Cell(20,2).Range.Text = "36+20=56"
$t.Cell(20,3).Range.Text = "59+33=92"
$t.Cell(20,4).Range.Text = "57+26=83"
$t.Cell(20,5).Range.Text = "36-26=10"
